$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 25; this shifts existing rows 25-54 down to 26-55
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new record
$ws.Cells.Item(25, 1).Value = 4
$ws.Cells.Item(25, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(25, 3).Value = "Los Lagos"
$ws.Cells.Item(25, 4).Value = 44579
$ws.Cells.Item(25, 5).Value = 10
$ws.Cells.Item(25, 6).Value = 100112031
$ws.Cells.Item(25, 7).Value = "Poroto verde"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 40
$ws.Cells.Item(25, 11).Value = 42000
$ws.Cells.Item(25, 12).Value = 42000
$ws.Cells.Item(25, 13).Value = 42000
$ws.Cells.Item(25, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(25, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(25, 16).Value = 1680
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"

# Match the date cell format used by the other date cells in column D
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
